$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename the header row: "<name>_old" -> "<name>_FV2404"
#    and "<name>_new" -> "<name>_FV2410" (columns A:J and L:U).
#    Column K keeps its header ("diff").
# ------------------------------------------------------------------
$cols = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $cols[$i] + "_FV2404"
}
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $cols[$i] + "_FV2410"
}

# ------------------------------------------------------------------
# 2) Turn the used range A1:U76 into an Excel Table ("Table1").
# ------------------------------------------------------------------
$tableRange = $ws.Range("A1:U76")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# ------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "edit applied"
